$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set each changed cell to Text format first so values such as "0.9981" or
# "26.734.28" are stored as literal strings (matching the original inlineStr
# cell type) instead of being auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.734.28'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.727.26'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9981'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.39'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9986'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4805'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2584'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06178'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.725.21'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '15.84'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06860'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6030'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.457'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.90'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9986'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.560.04'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.68%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007137'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.946.38'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.416'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.527'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.051'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.67'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.19'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.773'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.96%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '106.19'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.366'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.002'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07912'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.659'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04518'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.599'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9983'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6170'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.62%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9309'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.452'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.994'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.00%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01496'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.600'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.80'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3823'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.769'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1152'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05358'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.875'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.06'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.239'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.36'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.80%  '
